# Update the "想去人数" (column F) values in the "展览" and "全部类型"
# worksheets to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F on sheet "展览"
$sheet1Updates = @{
    3  = 16434
    5  = 736
    6  = 15613
    7  = 70
    8  = 9278
    9  = 493
    11 = 1034
    12 = 127
    13 = 225
    17 = 94
    18 = 619
    21 = 79
    22 = 1159
    23 = 15
    24 = 23
    26 = 536
    30 = 84
    32 = 66
    33 = 270
    34 = 369
    35 = 480
    36 = 122
    37 = 5708
    38 = 5250
}

# Row -> new value for column F on sheet "全部类型"
$sheet4Updates = @{
    3  = 16434
    5  = 736
    6  = 15613
    7  = 70
    8  = 9278
    9  = 493
    11 = 1034
    12 = 127
    13 = 225
    17 = 94
    18 = 619
    21 = 79
    22 = 1159
    23 = 15
    24 = 23
    26 = 536
    32 = 84
    34 = 66
    35 = 270
    36 = 369
    37 = 480
    38 = 122
    39 = 5708
    41 = 5250
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
